$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3176883.5
$ws.Range("I132").Value = 3923736.5
$ws.Range("J132").Value = 2757.8333
$ws.Range("K132").Value = 11771209.5
$ws.Range("L132").Value = 8273.499899999999
$ws.Range("M132").Value = -11768679.5
$ws.Range("N132").Value = -13333.4999
$ws.Range("H137").Value = 1564585.6
$ws.Range("I137").Value = 2176033.8
$ws.Range("J137").Value = 1996.1111
$ws.Range("K137").Value = 6528101.399999999
$ws.Range("L137").Value = 5988.3333
$ws.Range("M137").Value = -6525551.399999999
$ws.Range("N137").Value = -11088.3333
$ws.Range("H138").Value = 2575.3735
$ws.Range("I138").Value = 887.6383
$ws.Range("J138").Value = 4778.8057
$ws.Range("K138").Value = 2662.9149
$ws.Range("L138").Value = 14336.4171
$ws.Range("M138").Value = 2477.0851
$ws.Range("N138").Value = -24616.4171
$ws.Range("H141").Value = 218309.56
$ws.Range("I141").Value = 1078.1666
$ws.Range("J141").Value = 1738929.4
$ws.Range("K141").Value = 3234.4998
$ws.Range("L141").Value = 5216788.199999999
$ws.Range("M141").Value = 1945.5002
$ws.Range("N141").Value = -5227148.199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 22729744
$ws.Range("I2").Value = 35716310
$ws.Range("K2").Value = 35716310
$ws.Range("M2").Value = -35716197
$ws.Range("H44").Value = 20024.5
$ws.Range("J44").Value = 20024.5
$ws.Range("L44").Value = 20024.5
$ws.Range("N44").Value = -21000.5
$ws.Range("H61").Value = 1608.42
$ws.Range("I61").Value = 753.8919
$ws.Range("K61").Value = 753.8919
$ws.Range("M61").Value = -541.8919
$ws.Range("H74").Value = 842.85187
$ws.Range("I74").Value = 714
$ws.Range("J74").Value = 1409.8
$ws.Range("K74").Value = 714
$ws.Range("L74").Value = 1409.8
$ws.Range("M74").Value = 160
$ws.Range("N74").Value = -3157.8
$ws.Range("H77").Value = 842.85187
$ws.Range("I77").Value = 714
$ws.Range("J77").Value = 1409.8
$ws.Range("K77").Value = 3570
$ws.Range("L77").Value = 7049
$ws.Range("M77").Value = 798
$ws.Range("N77").Value = -15785
$ws.Range("H116").Value = 22729744
$ws.Range("I116").Value = 35716310
$ws.Range("K116").Value = 35716310
$ws.Range("M116").Value = -35714016
$ws.Range("H132").Value = 2324.8718
$ws.Range("I132").Value = 1762.5161
$ws.Range("J132").Value = 4504
$ws.Range("K132").Value = 5287.5483
$ws.Range("L132").Value = 13512
$ws.Range("M132").Value = -2757.5483
$ws.Range("N132").Value = -18572
$ws.Range("H136").Value = 1608.42
$ws.Range("I136").Value = 753.8919
$ws.Range("K136").Value = 2261.6757
$ws.Range("M136").Value = 288.3243000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 22729744
$ws.Range("I3").Value = 35716310
$ws.Range("K3").Value = 35716310
$ws.Range("M3").Value = -35716196
$ws.Range("H20").Value = 1483.8857
$ws.Range("I20").Value = 1338.909
$ws.Range("J20").Value = 1729.2307
$ws.Range("K20").Value = 1338.909
$ws.Range("L20").Value = 1729.2307
$ws.Range("M20").Value = -1091.909
$ws.Range("N20").Value = -2223.2307
$ws.Range("H134").Value = 1344.8108
$ws.Range("I134").Value = 617.9375
$ws.Range("J134").Value = 5996.8
$ws.Range("K134").Value = 1853.8125
$ws.Range("L134").Value = 17990.4
$ws.Range("M134").Value = 681.1875
$ws.Range("N134").Value = -23060.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2275876.2
$ws.Range("I31").Value = 3705699.8
$ws.Range("J31").Value = 4980.1177
$ws.Range("K31").Value = 3705699.8
$ws.Range("L31").Value = 4980.1177
$ws.Range("M31").Value = -3705404.8
$ws.Range("N31").Value = -5570.1177
$ws.Range("H34").Value = 2275876.2
$ws.Range("I34").Value = 3705699.8
$ws.Range("J34").Value = 4980.1177
$ws.Range("K34").Value = 3705699.8
$ws.Range("L34").Value = 4980.1177
$ws.Range("M34").Value = -3705497.8
$ws.Range("N34").Value = -5384.1177
$ws.Range("H58").Value = 8335456
$ws.Range("I58").Value = 1336.125
$ws.Range("J58").Value = 41671936
$ws.Range("K58").Value = 1336.125
$ws.Range("L58").Value = 41671936
$ws.Range("M58").Value = -1133.125
$ws.Range("N58").Value = -41672342
$ws.Range("H132").Value = 2265.1562
$ws.Range("I132").Value = 1618.8148
$ws.Range("K132").Value = 4856.4444
$ws.Range("M132").Value = -2326.4444
$ws.Range("H134").Value = 1694.9143
$ws.Range("I134").Value = 889.1852
$ws.Range("J134").Value = 4414.25
$ws.Range("K134").Value = 2667.5556
$ws.Range("L134").Value = 13242.75
$ws.Range("M134").Value = -132.5556000000001
$ws.Range("N134").Value = -18312.75
$ws.Range("H136").Value = 8335456
$ws.Range("I136").Value = 1336.125
$ws.Range("J136").Value = 41671936
$ws.Range("K136").Value = 4008.375
$ws.Range("L136").Value = 125015808
$ws.Range("M136").Value = -1458.375
$ws.Range("N136").Value = -125020908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 8898.571
$ws.Range("J34").Value = 10356.667
$ws.Range("L34").Value = 31070.001
$ws.Range("N34").Value = -31238.001
$ws.Range("H55").Value = 2898.48
$ws.Range("I55").Value = 380.8
$ws.Range("J55").Value = 3527.9
$ws.Range("K55").Value = 1142.4
$ws.Range("L55").Value = 10583.7
$ws.Range("M55").Value = -965.4000000000001
$ws.Range("N55").Value = -10937.7
$ws.Range("H130").Value = 2338.3333
$ws.Range("J130").Value = 2600
$ws.Range("L130").Value = 7800
$ws.Range("N130").Value = -17840
$ws.Range("H140").Value = 11908962
$ws.Range("I140").Value = 33334114
$ws.Range("J140").Value = 6100
$ws.Range("K140").Value = 100002342
$ws.Range("L140").Value = 18300
$ws.Range("M140").Value = -99997162
$ws.Range("N140").Value = -28660

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2778.138
$ws.Range("I132").Value = 2369.8096
$ws.Range("J132").Value = 3850
$ws.Range("K132").Value = 7109.4288
$ws.Range("L132").Value = 11550
$ws.Range("M132").Value = -4579.4288
$ws.Range("N132").Value = -16610
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 962.26666
$ws.Range("I55").Value = 176.66667
$ws.Range("J55").Value = 1486
$ws.Range("K55").Value = 176.66667
$ws.Range("L55").Value = 1486
$ws.Range("M55").Value = -3.666670000000011
$ws.Range("N55").Value = -1832
$ws.Range("H93").Value = 2478.1904
$ws.Range("I93").Value = 2131.7058
$ws.Range("J93").Value = 3950.75
$ws.Range("K93").Value = 2131.7058
$ws.Range("L93").Value = 3950.75
$ws.Range("M93").Value = -883.7058000000002
$ws.Range("N93").Value = -6446.75
$ws.Range("H122").Value = 3185.1853
$ws.Range("I122").Value = 2761.9048
$ws.Range("K122").Value = 8285.714399999999
$ws.Range("M122").Value = -5835.714399999999
$ws.Range("H132").Value = 2765.0908
$ws.Range("I132").Value = 1788.8
$ws.Range("K132").Value = 5366.4
$ws.Range("M132").Value = -2836.4
$ws.Range("H136").Value = 2705675.2
$ws.Range("I136").Value = 3573144.8
$ws.Range("J136").Value = 6881.6665
$ws.Range("K136").Value = 10719434.4
$ws.Range("L136").Value = 20644.9995
$ws.Range("M136").Value = -10716884.4
$ws.Range("N136").Value = -25744.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 11898.182
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 11898.182
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 11898.182
$ws.Range("N54").Value = -12938.182
$ws.Range("M54").ClearContents()
$ws.Range("H132").Value = 252774
$ws.Range("I132").Value = 373538.88
$ws.Range("K132").Value = 1120616.64
$ws.Range("M132").Value = -1118086.64
$ws.Range("H136").Value = 924.72546
$ws.Range("I136").Value = 509.85184
$ws.Range("J136").Value = 1391.4584
$ws.Range("K136").Value = 1529.55552
$ws.Range("L136").Value = 4174.3752
$ws.Range("M136").Value = 1020.44448
$ws.Range("N136").Value = -9274.3752
